$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (C) column date for rows 2 through 19 from 2023-09-06 (45175)
# to 2023-09-14 (45183), keeping the existing date formatting/style intact.
for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
